$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new row 16: LeetCode 2014 "Longest Subsequence Repeated k Times" ---
$ws.Cells.Item(16, 1).Value = 2014
$ws.Cells.Item(16, 2).Value = "Longest Subsequence Repeated k Times"
$ws.Cells.Item(16, 3).Value = "#string #backtracking #greedy #enumeration "
$ws.Cells.Item(16, 4).Value = "hard"
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 56
$ws.Cells.Item(16, 8).Value = 45835
$ws.Cells.Item(16, 9).Value = 45835

# Copy the date formatting (style) from the row above for H16:I16 so the new
# "First"/"Last Update" cells render as dates instead of raw serials, without
# introducing a brand-new (duplicate) number-format style.
$ws.Range("H15:I15").Copy()
$ws.Range("H16:I16").PasteSpecial(-4122)

# Row 16 needs to be tall enough to show the wrapped 4-line tag text in C16.
$ws.Rows(16).RowHeight = 68

# --- View state: move the selection/scroll position down to the new row ---
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I16").Select()
